# Update TPM-derived values on Sheet1 with newly computed results.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 1.918906333333333
$ws.Range("N2").Value = 5.756718999999999
$ws.Range("O2").Value = 0.006524019162508824
$ws.Range("P2").Value = 0.006524019162508824
$ws.Range("Q2").Value = 0.1250973416826666
$ws.Range("R2").Value = 1.125876075144
$ws.Range("S2").Value = 0.006524019162508824
$ws.Range("T2").Value = 0.006524019162508824

# Row 3
$ws.Range("O3").Value = 0.6163557430885885
$ws.Range("P3").Value = 0.6163557430885885
$ws.Range("S3").Value = 0.6163557430885885
$ws.Range("T3").Value = 0.6163557430885885

# Row 4
$ws.Range("M4").Value = 29.04767233333333
$ws.Range("N4").Value = 87.143017
$ws.Range("O4").Value = 0.09875811426384234
$ws.Range("P4").Value = 0.09875811426384236
$ws.Range("Q4").Value = 1.893675854754667
$ws.Range("R4").Value = 17.043082692792
$ws.Range("S4").Value = 0.09875811426384234
$ws.Range("T4").Value = 0.09875811426384236

# Row 5
$ws.Range("M5").Value = 81.87450533333333
$ws.Range("N5").Value = 245.623516
$ws.Range("O5").Value = 0.2783621234850603
$ws.Range("P5").Value = 0.2783621234850603
$ws.Range("Q5").Value = 5.337562751690666
$ws.Range("R5").Value = 48.038064765216
$ws.Range("S5").Value = 0.2783621234850603
$ws.Range("T5").Value = 0.2783621234850603
